# Scheduled-runner refresh of Leve vendor/market price & profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the ALC, ARM, BSM,
# CRP, CUL and GSM crafter leve tables. Values below are plain numbers
# (no formulas are used in this workbook), so each target cell is simply
# re-written with its refreshed figure; a few rows also gained/lost a
# trailing LeveProfitHQ (N) cell because HQ pricing became known/unknown.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 18
$ws.Range("H18").Value = 1000000
$ws.Range("I18").Value = 1000000
$ws.Range("K18").Value = 1000000
$ws.Range("M18").Value = -999716
# row 19
$ws.Range("H19").Value = 1709.2
$ws.Range("I19").Value = 1739.8182
$ws.Range("J19").Value = 1625
$ws.Range("K19").Value = 1739.8182
$ws.Range("L19").Value = 1625
$ws.Range("M19").Value = -1564.8182
$ws.Range("N19").Value = -1975
# row 32
$ws.Range("H32").Value = 9699.799999999999
$ws.Range("J32").Value = 10499.667
$ws.Range("L32").Value = 10499.667
$ws.Range("N32").Value = -11151.667
# row 33
$ws.Range("H33").Value = 126.9
$ws.Range("I33").Value = 126.9
$ws.Range("K33").Value = 126.9
$ws.Range("M33").Value = 102.1
# row 39
$ws.Range("H39").Value = 279.8
$ws.Range("I39").Value = 279.8
$ws.Range("K39").Value = 839.4000000000001
$ws.Range("M39").Value = -543.4000000000001
# row 92
$ws.Range("H92").Value = 384.3846
$ws.Range("I92").Value = 366.8889
$ws.Range("K92").Value = 366.8889
$ws.Range("M92").Value = 881.1111000000001
# row 106
$ws.Range("H106").Value = 10000
$ws.Range("I106").Value = 10000
$ws.Range("K106").Value = 10000
$ws.Range("M106").Value = -9369
# row 111
$ws.Range("H111").Value = 8333
$ws.Range("I111").Value = 8499.5
$ws.Range("J111").Value = 8000
$ws.Range("K111").Value = 25498.5
$ws.Range("L111").Value = 24000
$ws.Range("M111").Value = -22431.5
$ws.Range("N111").Value = -30134
# row 132
$ws.Range("H132").Value = 2979.55
$ws.Range("I132").Value = 2353.7693
$ws.Range("K132").Value = 7061.3079
$ws.Range("M132").Value = -4531.3079
# row 137
$ws.Range("H137").Value = 1596.4615
$ws.Range("I137").Value = 1347.75
$ws.Range("K137").Value = 4043.25
$ws.Range("M137").Value = -1493.25
# row 138
$ws.Range("H138").Value = 2700.641
$ws.Range("J138").Value = 3057.4443
$ws.Range("L138").Value = 9172.332900000001
$ws.Range("N138").Value = -19452.3329

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 6694.5
$ws.Range("I32").Value = 4920.467
$ws.Range("K32").Value = 4920.467
$ws.Range("M32").Value = -4633.467
# row 74
$ws.Range("H74").Value = 33323698
$ws.Range("I74").Value = 49982300
$ws.Range("K74").Value = 49982300
$ws.Range("M74").Value = -49981426
# row 77
$ws.Range("H77").Value = 33323698
$ws.Range("I77").Value = 49982300
$ws.Range("K77").Value = 249911500
$ws.Range("M77").Value = -249907132
# row 102
$ws.Range("H102").Value = 1566.1666
$ws.Range("I102").Value = 1566.1666
$ws.Range("K102").Value = 1566.1666
$ws.Range("M102").Value = 55.83339999999998

$ws = $wb.Worksheets.Item("BSM")
# row 22
$ws.Range("H22").Value = 669.53845
$ws.Range("I22").Value = 636.9091
$ws.Range("K22").Value = 636.9091
$ws.Range("M22").Value = -463.9091
# row 86
$ws.Range("H86").Value = 5360.778
$ws.Range("I86").Value = 3041.3333
$ws.Range("J86").Value = 9999.666999999999
$ws.Range("K86").Value = 3041.3333
$ws.Range("L86").Value = 9999.666999999999
$ws.Range("M86").Value = -1918.3333
$ws.Range("N86").Value = -12245.667
# row 89
$ws.Range("H89").Value = 5360.778
$ws.Range("I89").Value = 3041.3333
$ws.Range("J89").Value = 9999.666999999999
$ws.Range("K89").Value = 15206.6665
$ws.Range("L89").Value = 49998.335
$ws.Range("M89").Value = -9590.666499999999
$ws.Range("N89").Value = -61230.335
# row 94
$ws.Range("H94").Value = 406.91666
$ws.Range("I94").Value = 443.9
$ws.Range("K94").Value = 443.9
$ws.Range("M94").Value = 7.100000000000023
# row 105
$ws.Range("H105").Value = 3167.0833
$ws.Range("I105").Value = 2881.1
$ws.Range("K105").Value = 2881.1
$ws.Range("M105").Value = -1134.1
# row 135
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140

$ws = $wb.Worksheets.Item("CRP")
# row 2
$ws.Range("H2").Value = 2575.8
$ws.Range("J2").Value = 3999.6667
$ws.Range("L2").Value = 3999.6667
$ws.Range("N2").Value = -4225.6667
# row 5
$ws.Range("H5").Value = 997.5
$ws.Range("J5").Value = 997.5
$ws.Range("L5").Value = 997.5
$ws.Range("N5").Value = -1221.5
# row 6
$ws.Range("H6").Value = 1000
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
# row 7
$ws.Range("H7").Value = 474.16666
$ws.Range("I7").Value = 286.25
$ws.Range("K7").Value = 286.25
$ws.Range("M7").Value = -173.25
# row 12
$ws.Range("H12").Value = 3000
$ws.Range("J12").Value = 3000
$ws.Range("L12").Value = 3000
$ws.Range("N12").Value = -3340
# row 17
$ws.Range("H17").Value = 16336
$ws.Range("J17").Value = 17004.5
$ws.Range("L17").Value = 17004.5
$ws.Range("N17").Value = -17352.5
# row 25
$ws.Range("H25").Value = 20006.5
$ws.Range("J25").Value = 20006.5
$ws.Range("L25").Value = 20006.5
$ws.Range("N25").Value = -20354.5
# row 41
$ws.Range("H41").Value = 39965
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 39965
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 39965
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -40821
# row 51
$ws.Range("H51").Value = 50099
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 50099
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 50099
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -51571
# row 60
$ws.Range("H60").Value = 35328.12
$ws.Range("I60").Value = 18000
$ws.Range("J60").Value = 36834.914
$ws.Range("K60").Value = 18000
$ws.Range("L60").Value = 36834.914
$ws.Range("M60").Value = -17489
$ws.Range("N60").Value = -37856.914
# row 61
$ws.Range("H61").Value = 50099
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 50099
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 50099
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -50795
# row 99
$ws.Range("H99").Value = 5999
$ws.Range("J99").Value = 6000
$ws.Range("L99").Value = 6000
$ws.Range("N99").Value = -8996
# row 126
$ws.Range("H126").Value = 5999
$ws.Range("J126").Value = 6000
$ws.Range("L126").Value = 18000
$ws.Range("N126").Value = -22940
# row 134
$ws.Range("H134").Value = 3054.1052
$ws.Range("I134").Value = 2931.8572
$ws.Range("K134").Value = 8795.571599999999
$ws.Range("M134").Value = -6260.571599999999
# row 141
$ws.Range("H141").Value = 60000
$ws.Range("J141").Value = 60000
$ws.Range("L141").Value = 60000
$ws.Range("N141").Value = -70360

$ws = $wb.Worksheets.Item("CUL")
# row 56
$ws.Range("H56").Value = 18416.883
$ws.Range("I56").Value = 18416.883
$ws.Range("K56").Value = 18416.883
$ws.Range("M56").Value = -17886.883
# row 122
$ws.Range("H122").Value = 101641.1
$ws.Range("I122").Value = 1496
$ws.Range("J122").Value = 201786.2
$ws.Range("K122").Value = 13464
$ws.Range("L122").Value = 1816075.8
$ws.Range("M122").Value = -11014
$ws.Range("N122").Value = -1820975.8

$ws = $wb.Worksheets.Item("GSM")
# row 2
$ws.Range("H2").Value = 92.55556
$ws.Range("I2").Value = 87.71429000000001
$ws.Range("J2").Value = 109.5
$ws.Range("K2").Value = 87.71429000000001
$ws.Range("L2").Value = 109.5
$ws.Range("M2").Value = 25.28570999999999
$ws.Range("N2").Value = -335.5
# row 46
$ws.Range("H46").Value = 34941.668
$ws.Range("J46").Value = 34941.668
$ws.Range("L46").Value = 34941.668
$ws.Range("N46").Value = -35253.668
# row 57
$ws.Range("H57").Value = 26666
$ws.Range("I57").Value = 3998
$ws.Range("K57").Value = 3998
$ws.Range("M57").Value = -3178
# row 80
$ws.Range("H80").Value = 5947.8
$ws.Range("I80").Value = 5328.6665
$ws.Range("J80").Value = 6876.5
$ws.Range("K80").Value = 5328.6665
$ws.Range("L80").Value = 6876.5
$ws.Range("M80").Value = -4330.6665
$ws.Range("N80").Value = -8872.5
# row 83
$ws.Range("H83").Value = 5947.8
$ws.Range("I83").Value = 5328.6665
$ws.Range("J83").Value = 6876.5
$ws.Range("K83").Value = 26643.3325
$ws.Range("L83").Value = 34382.5
$ws.Range("M83").Value = -21651.3325
$ws.Range("N83").Value = -44366.5
# row 102
$ws.Range("H102").Value = 1121.95
$ws.Range("I102").Value = 1144.3334
$ws.Range("K102").Value = 1144.3334
$ws.Range("M102").Value = 477.6666
